# Generate Report for Handback
# Updates the localization-status report: marks rows as handed back,
# records the handback datetime/target-file/handback-file for each
# language, and widens a few columns to fit the new (longer) content.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

$hoMd405 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4b6a43944c4685eddf618138537d2a02110d58ba/e2e/405fef66-4081-465f-9681-c573232c0bf6.md"
$hoMdAa0 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4b6a43944c4685eddf618138537d2a02110d58ba/e2e/aa0d2dcb-005b-4116-9c5f-da4b2a48f3fe.md"

# ---------------------------------------------------------------------
# 1. Status rolls from "In Translation" to "Handed back: in sync with
#    en-US" everywhere it is shown: the Overview rollup (per language)
#    and each language detail sheet's Status column.
# ---------------------------------------------------------------------
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn detail sheet: fill in Latest Target File / Latest Handback
#    File / Latest Handback DateTime for both rows, and add hyperlinks
#    on the new "Latest Target File" cells (same link target as the
#    Source File Name column).
# ---------------------------------------------------------------------
$wsZh.Range("I2").Value = "405fef66-4081-465f-9681-c573232c0bf6.md"
$wsZh.Range("J2").Value = "405fef66-4081-465f-9681-c573232c0bf6.8ec233ee61fc41b48961d4acbafbf2dae9a32a7f.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-24 14:26:34"

$wsZh.Range("I3").Value = "aa0d2dcb-005b-4116-9c5f-da4b2a48f3fe.md"
$wsZh.Range("J3").Value = "aa0d2dcb-005b-4116-9c5f-da4b2a48f3fe.7637a1d2c843eee9867237d2f2fea00f39266d41.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-24 14:26:34"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hoMd405, "", "", "405fef66-4081-465f-9681-c573232c0bf6.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $hoMd405, "", "", "405fef66-4081-465f-9681-c573232c0bf6.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $hoMdAa0, "", "", "aa0d2dcb-005b-4116-9c5f-da4b2a48f3fe.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $hoMdAa0, "", "", "aa0d2dcb-005b-4116-9c5f-da4b2a48f3fe.md")

# ---------------------------------------------------------------------
# 3. de-de detail sheet: same shape of update, with de-de handoff/
#    handback file names and its own handback timestamp.
# ---------------------------------------------------------------------
$wsDe.Range("I2").Value = "405fef66-4081-465f-9681-c573232c0bf6.md"
$wsDe.Range("J2").Value = "405fef66-4081-465f-9681-c573232c0bf6.8ec233ee61fc41b48961d4acbafbf2dae9a32a7f.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-24 14:26:41"

$wsDe.Range("I3").Value = "aa0d2dcb-005b-4116-9c5f-da4b2a48f3fe.md"
$wsDe.Range("J3").Value = "aa0d2dcb-005b-4116-9c5f-da4b2a48f3fe.7637a1d2c843eee9867237d2f2fea00f39266d41.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-24 14:26:41"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hoMd405, "", "", "405fef66-4081-465f-9681-c573232c0bf6.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $hoMd405, "", "", "405fef66-4081-465f-9681-c573232c0bf6.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $hoMdAa0, "", "", "aa0d2dcb-005b-4116-9c5f-da4b2a48f3fe.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $hoMdAa0, "", "", "aa0d2dcb-005b-4116-9c5f-da4b2a48f3fe.md")

# ---------------------------------------------------------------------
# 4. Widen columns that now hold longer content.
# ---------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 29.14
$ws1.Columns.Item(6).ColumnWidth = 29.14

$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
